# add prodi & kuota dosen
# Adds a new "program_studi" column (H) with a repeating TRPL/TRK/BD
# pattern for the dosen data rows, and updates the sheet view (zoom /
# selection) plus column width to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H: header + data -------------------------------------
$ws.Range("H1").Value = "program_studi"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").VerticalAlignment = -4108
$ws.Range("H1").HorizontalAlignment = -4108

$pattern = @("TRPL", "TRK", "BD")

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.Value = $pattern[($row - 2) % 3]
}

# Style the first data cell directly, then fan the same formatting out to
# the rest of the column via copy/paste-special so every cell collapses
# onto a single shared style (matches s="19" for H2:H23 in the target).
$ws.Range("H2").VerticalAlignment = -4108
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").Copy()
$ws.Range("H3:H23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column width for H ------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 18.25

# --- Sheet view: zoom + selection ---------------------------------------
$excel.ActiveWindow.Zoom = 87
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K9").Select()
